$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix up the "5/2/2023" -> "5/02/2023" date text entries and rotate the
# assignment notes between rows 5, 6, 7 and 20.
$ws.Range("C6").Value = "Tuấn:5/02/2023"
$ws.Range("C5").Value = "Dũng:5/02/2023"
$ws.Range("C7").Value = "Nam:5/02/2023"
$ws.Range("C20").Value = "Nam:31/01/2023"

# B8 becomes a numeric 100% value instead of the placeholder "?%" text.
$ws.Range("B8").Value = 1
$ws.Range("B8").NumberFormat = "0%"

# Update the active selection to A24.
$ws.Range("A24").Select()
